# river update May 2024
# Appends 4 new result rows (36-39) for "Tokiahuru at Karioi Domain Road"
# covering the 2019 - 2023 period: ASPM, Chlorophyll A (92nd Percentile),
# MCI and QMCI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common columns shared by every appended row.
$siteName = "Tokiahuru at Karioi Domain Road"
$yearRange = "2019 - 2023"
$siteType = "RepSite"
$nztmX = 1815033.57
$nztmY = 5627502.8
$district = "Ruapehu District"
$fmu = "Whangaehu"
$zone = "Upper Whangaehu"
$subzone = "Whau_1c"

# --- Row 36: ASPM ---
$ws.Cells.Item(36, 1).Value = $siteName
$ws.Cells.Item(36, 2).Value = "ASPM"
$ws.Cells.Item(36, 3).Value = "C"
$ws.Cells.Item(36, 4).Value = $yearRange
$ws.Cells.Item(36, 5).Value = $siteType
$ws.Cells.Item(36, 6).Value = 0.37
$ws.Cells.Item(36, 7).Value = 0.3546
$ws.Cells.Item(36, 8).Value = 0.42
$ws.Cells.Item(36, 9).Value = 0.42
$ws.Cells.Item(36, 12).Value = 0.3565
$ws.Cells.Item(36, 13).Value = 0.40915
$ws.Cells.Item(36, 14).Value = 0.42
$ws.Cells.Item(36, 15).Value = $nztmX
$ws.Cells.Item(36, 16).Value = $nztmY
$ws.Cells.Item(36, 17).Value = $district
$ws.Cells.Item(36, 18).Value = $fmu
$ws.Cells.Item(36, 19).Value = $zone
$ws.Cells.Item(36, 20).Value = $subzone

# --- Row 37: Chlorophyll A (92nd Percentile) ---
$ws.Cells.Item(37, 1).Value = $siteName
$ws.Cells.Item(37, 2).Value = "Chlorophyll A (92nd Percentile)"
$ws.Cells.Item(37, 3).Value = "B"
$ws.Cells.Item(37, 4).Value = $yearRange
$ws.Cells.Item(37, 5).Value = $siteType
$ws.Cells.Item(37, 6).Value = 19.5
$ws.Cells.Item(37, 7).Value = 29.6562244897959
$ws.Cells.Item(37, 8).Value = 165
$ws.Cells.Item(37, 9).Value = 70
$ws.Cells.Item(37, 12).Value = 21.25
$ws.Cells.Item(37, 13).Value = 55.85
$ws.Cells.Item(37, 14).Value = 67.90000000000001
$ws.Cells.Item(37, 15).Value = $nztmX
$ws.Cells.Item(37, 16).Value = $nztmY
$ws.Cells.Item(37, 17).Value = $district
$ws.Cells.Item(37, 18).Value = $fmu
$ws.Cells.Item(37, 19).Value = $zone
$ws.Cells.Item(37, 20).Value = $subzone
$ws.Cells.Item(37, 21).Value = "mg chl-a /m2"

# --- Row 38: MCI ---
$ws.Cells.Item(38, 1).Value = $siteName
$ws.Cells.Item(38, 2).Value = "MCI"
$ws.Cells.Item(38, 3).Value = "B"
$ws.Cells.Item(38, 4).Value = $yearRange
$ws.Cells.Item(38, 5).Value = $siteType
$ws.Cells.Item(38, 6).Value = 111
$ws.Cells.Item(38, 7).Value = 104.968
$ws.Cells.Item(38, 8).Value = 114.4
$ws.Cells.Item(38, 9).Value = 114.4
$ws.Cells.Item(38, 12).Value = 111.5
$ws.Cells.Item(38, 13).Value = 113.56
$ws.Cells.Item(38, 14).Value = 114.4
$ws.Cells.Item(38, 15).Value = $nztmX
$ws.Cells.Item(38, 16).Value = $nztmY
$ws.Cells.Item(38, 17).Value = $district
$ws.Cells.Item(38, 18).Value = $fmu
$ws.Cells.Item(38, 19).Value = $zone
$ws.Cells.Item(38, 20).Value = $subzone

# --- Row 39: QMCI ---
$ws.Cells.Item(39, 1).Value = $siteName
$ws.Cells.Item(39, 2).Value = "QMCI"
$ws.Cells.Item(39, 3).Value = "D"
$ws.Cells.Item(39, 4).Value = $yearRange
$ws.Cells.Item(39, 5).Value = $siteType
$ws.Cells.Item(39, 6).Value = 4
$ws.Cells.Item(39, 7).Value = 4.0338
$ws.Cells.Item(39, 8).Value = 5.18
$ws.Cells.Item(39, 9).Value = 5.18
$ws.Cells.Item(39, 12).Value = 3.8
$ws.Cells.Item(39, 13).Value = 4.886
$ws.Cells.Item(39, 14).Value = 5.18
$ws.Cells.Item(39, 15).Value = $nztmX
$ws.Cells.Item(39, 16).Value = $nztmY
$ws.Cells.Item(39, 17).Value = $district
$ws.Cells.Item(39, 18).Value = $fmu
$ws.Cells.Item(39, 19).Value = $zone
$ws.Cells.Item(39, 20).Value = $subzone
